$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "36.353.24"
$ws.Cells.Item(2, 5).Value = "  +0.14%  "

$ws.Cells.Item(3, 4).Value = "1.934.28"
$ws.Cells.Item(3, 5).Value = "  -2.05%  "

$ws.Cells.Item(4, 5).Value = "  -0.15%  "

$ws.Cells.Item(5, 4).Value = "'241.63"
$ws.Cells.Item(5, 5).Value = "  -1.27%  "

$ws.Cells.Item(6, 4).Value = "'0.607"
$ws.Cells.Item(6, 5).Value = "  -2.56%  "

$ws.Cells.Item(7, 5).Value = "  -0.08%  "

$ws.Cells.Item(8, 4).Value = "'56.37"
$ws.Cells.Item(8, 5).Value = "  -3.62%  "

$ws.Cells.Item(9, 5).Value = "  -3.50%  "

$ws.Cells.Item(10, 4).Value = "'0.0838"
$ws.Cells.Item(10, 5).Value = "  +0.67%  "

$ws.Cells.Item(11, 5).Value = "  -1.65%  "

$ws.Cells.Item(12, 4).Value = "2.217.85"
$ws.Cells.Item(12, 5).Value = "  -2.07%  "

$ws.Cells.Item(13, 4).Value = "'21.09"
$ws.Cells.Item(13, 5).Value = "  -7.93%  "

$ws.Cells.Item(14, 4).Value = "'0.799"
$ws.Cells.Item(14, 5).Value = "  -6.39%  "

$ws.Cells.Item(15, 4).Value = "'13.34"
$ws.Cells.Item(15, 5).Value = "  -3.60%  "

$ws.Cells.Item(16, 5).Value = "  -5.38%  "

$ws.Cells.Item(17, 4).Value = "1.933.09"
$ws.Cells.Item(17, 5).Value = "  -3.25%  "

$ws.Cells.Item(18, 4).Value = "36.291.37"
$ws.Cells.Item(18, 5).Value = "  +0.26%  "

$ws.Cells.Item(19, 4).Value = "0.0₃0859"
$ws.Cells.Item(19, 5).Value = "  -2.17%  "

$ws.Cells.Item(20, 4).Value = "'68.70"
$ws.Cells.Item(20, 5).Value = "  -2.17%  "

$ws.Cells.Item(21, 4).Value = "'226.05"
$ws.Cells.Item(21, 5).Value = "  -2.97%  "

$ws.Cells.Item(22, 4).Value = "'4.93"
$ws.Cells.Item(22, 5).Value = "  -5.88%  "

$ws.Cells.Item(23, 5).Value = "  -0.17%  "

$ws.Cells.Item(24, 5).Value = "  -6.70%  "

$ws.Cells.Item(25, 4).Value = "'2.26"
$ws.Cells.Item(25, 5).Value = "  -3.32%  "

$ws.Cells.Item(26, 4).Value = "'9.06"
$ws.Cells.Item(26, 5).Value = "  -7.09%  "

$ws.Cells.Item(27, 4).Value = "'160.76"
$ws.Cells.Item(27, 5).Value = "  -1.45%  "

$ws.Cells.Item(28, 4).Value = "'0.131"
$ws.Cells.Item(28, 5).Value = "  -1.95%  "

$ws.Cells.Item(29, 5).Value = "  -3.09%  "

$ws.Cells.Item(30, 5).Value = "  -1.79%  "

$ws.Cells.Item(31, 5).Value = "  -6.21%  "

$ws.Cells.Item(32, 4).Value = "'4.51"
$ws.Cells.Item(32, 5).Value = "  -6.82%  "

$ws.Cells.Item(33, 4).Value = "'0.0616"
$ws.Cells.Item(33, 5).Value = "  -8.13%  "

$ws.Cells.Item(34, 5).Value = "  -5.59%  "

$ws.Cells.Item(35, 5).Value = "  -0.01%  "

$ws.Cells.Item(36, 4).Value = "'5.92"
$ws.Cells.Item(36, 5).Value = "  -2.82%  "

$ws.Cells.Item(37, 4).Value = "'1.77"
$ws.Cells.Item(37, 5).Value = "  -2.09%  "

$ws.Cells.Item(38, 5).Value = "  -3.61%  "

$ws.Cells.Item(39, 4).Value = "'2.95"
$ws.Cells.Item(39, 5).Value = "  +2.17%  "

$ws.Cells.Item(40, 4).Value = "'0.0967"
$ws.Cells.Item(40, 5).Value = "  +1.11%  "

$ws.Cells.Item(41, 4).Value = "'2.87"
$ws.Cells.Item(41, 5).Value = "  -0.74%  "

$ws.Cells.Item(42, 4).Value = "'0.0208"
$ws.Cells.Item(42, 5).Value = "  -2.09%  "

$ws.Cells.Item(43, 5).Value = "  -6.27%  "

$ws.Cells.Item(44, 4).Value = "'15.51"
$ws.Cells.Item(44, 5).Value = "  -2.91%  "

$ws.Cells.Item(45, 4).Value = "1.323.23"
$ws.Cells.Item(45, 5).Value = "  -2.49%  "

$ws.Cells.Item(46, 5).Value = "  -6.44%  "

$ws.Cells.Item(47, 4).Value = "'84.81"
$ws.Cells.Item(47, 5).Value = "  -6.96%  "

$ws.Cells.Item(48, 5).Value = "  -5.08%  "

$ws.Cells.Item(50, 4).Value = "2.108.72"
$ws.Cells.Item(50, 5).Value = "  -2.05%  "

$ws.Cells.Item(51, 4).Value = "'43.18"
$ws.Cells.Item(51, 5).Value = "  -3.12%  "
